$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("R2").Value = -506.7265263990728
$ws.Range("S2").Value = -36.26819966954229
$ws.Range("T2").Value = -0.9698828417572424
$ws.Range("U2").Value = 0.02509995865046097
$ws.Range("V2").Value = -0.04766639009244159
$ws.Range("W2").Value = 0.005455193700233254
$ws.Range("X2").Value = -0.02365414841562918
$ws.Range("Y2").Value = -0.01860716671230337
$ws.Range("Z2").Value = -0.04074048993808185
$ws.Range("AA2").Value = 0.08111972467532483
$ws.Range("AB2").Value = 0.005526977671832082
$ws.Range("AC2").Value = 0.0160868691478208
$ws.Range("AD2").Value = -0.07235956285197402
$ws.Range("AE2").Value = 0.006584135528506257
$ws.Range("AF2").Value = -0.1011787126675776

# Row 3
$ws.Range("R3").Value = -227.365185202805
$ws.Range("S3").Value = -19.74355796644368
$ws.Range("T3").Value = -3.931883658503619
$ws.Range("U3").Value = 0.1668804729245199
$ws.Range("V3").Value = 0.0507993524914901
$ws.Range("W3").Value = -0.007036138791104213
$ws.Range("X3").Value = -0.005296282186903855
$ws.Range("Y3").Value = -0.01230464106693177
$ws.Range("Z3").Value = 0.003152918567896873
$ws.Range("AA3").Value = 0.02145437210630252
$ws.Range("AB3").Value = 0.002411604044184461
$ws.Range("AC3").Value = -0.001497776909519523
$ws.Range("AD3").Value = -0.004795516862897471
$ws.Range("AE3").Value = 0.006386182358791445
$ws.Range("AF3").Value = 0.0005124172718750669

# Row 4
$ws.Range("R4").Value = 0.6543114669707606
$ws.Range("S4").Value = 0.01841955344801832
$ws.Range("T4").Value = 0.04997780741549365
$ws.Range("U4").Value = -0.0003998127625309534
$ws.Range("V4").Value = -0.005685289958183402
$ws.Range("W4").Value = -0.0001588694580014182
$ws.Range("X4").Value = -0.003134602381178558
$ws.Range("Y4").Value = -0.000004995518037629515
$ws.Range("Z4").Value = 0.002854967330332574
$ws.Range("AA4").Value = 0.000004390874723644253
$ws.Range("AB4").Value = -0.001539320870255615
$ws.Range("AC4").Value = 0.0000008104266509614035
$ws.Range("AD4").Value = 0.0006066464535303284
$ws.Range("AE4").Value = -0.000002082798250929089
$ws.Range("AF4").Value = -0.0001603569647642223

# Row 5
$ws.Range("R5").Value = 105.2800992812356
$ws.Range("S5").Value = 0.7832715440462812
$ws.Range("T5").Value = -8.240258121134568
$ws.Range("U5").Value = -0.004149490925661859
$ws.Range("V5").Value = -3.738582750049267
$ws.Range("W5").Value = -0.003592570794600292
$ws.Range("X5").Value = -2.637583783252583
$ws.Range("Y5").Value = 0.001450863138335099
$ws.Range("Z5").Value = -1.385926408769692
$ws.Range("AA5").Value = -0.0009994469405432607
$ws.Range("AB5").Value = -0.8785102085623572
$ws.Range("AC5").Value = -0.0003560227665888279
$ws.Range("AD5").Value = -0.353241166349535
$ws.Range("AE5").Value = -0.0002076459765332894
$ws.Range("AF5").Value = -0.1917896709063655

# Row 6
$ws.Range("R6").Value = 269.7823546950096
$ws.Range("S6").Value = 3.855652840467855
$ws.Range("T6").Value = -0.5729623834209503
$ws.Range("U6").Value = -0.0008179500907845318
$ws.Range("V6").Value = -0.02233037482404465
$ws.Range("W6").Value = -0.01318264973593808
$ws.Range("X6").Value = -0.0137540298274969
$ws.Range("Y6").Value = -0.00006425767928124116
$ws.Range("Z6").Value = -0.008907175013527993
$ws.Range("AA6").Value = 0.01827174553327056
$ws.Range("AB6").Value = -0.02700939565639115
$ws.Range("AC6").Value = -0.007658160588164015
$ws.Range("AD6").Value = -0.024789631445696
$ws.Range("AE6").Value = -0.01352360735764606
$ws.Range("AF6").Value = -0.002060793399062003
